$d = $word.ActiveDocument

# 1. Merge the three runs that together spell out
#    "{{ surname }} {{ name }} {{ middle_name }}" into a single run of text.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("{{ surname }} {{ name }} {{ middle_name }}", $true, $false, $false, $false, $false, `
              $true, 1, $false, "{{ surname }} {{ name }} {{ middle_name }}", 2)

# 2. Lengthen the run of spaces after "ИУЦТ" in the "отдела ИУЦТ" line so it
#    grows from 76 to 79 trailing spaces (3 extra spaces), to line up with the
#    other director/representative name lines.
$old = "отдела ИУЦТ" + "".PadRight(76)
$new = "отдела ИУЦТ" + "".PadRight(79)
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute($old, $true, $false, $false, $false, $false, `
               $true, 1, $false, $new, 2)
